# Update "Check Ve TBD" test case sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new test-case row values in the same order the strings were
# originally authored so the shared-string table indices line up:
# B12, A12, C12, D12, E12, F12, G12
$ws.Range("B12").Value = "Kiểm tra link Về TDB"
$ws.Range("A12").Value = "Check Về TBD"
$ws.Range("C12").Value = "baseURL: https://tbd.edu.vn/"
$ws.Range("D12").Value = "Text của h1.titleMainPage"
$ws.Range("E12").Value = "Về TBD"
$ws.Range("F12").Value = "Pass"
$ws.Range("G12").Value = "Đức"

# Adjust column widths (closest values the width->pixel-grid rounding of
# this COM layer maps onto the target OOXML widths of 27.42578125 / 14.85546875)
$ws.Range("C1:D1").ColumnWidth = 26.666666666666668
$ws.Range("G1").ColumnWidth = 14.0

# Update selected cell
$ws.Range("C11").Select()
